$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append after the existing last row (229).
# Columns: A = date serial (style carried over from the rows above),
#          B = nuovi pos., C = somma mobile 7gg., D = somma mobile 7gg. per 100mila abitanti
$rows = @(
    @{ Row = 230; A = 44304; B = 1; C = 4; D = 87.24100327153762 },
    @{ Row = 231; A = 44305; B = 0; C = 2; D = 43.62050163576881 },
    @{ Row = 232; A = 44306; B = 0; C = 2; D = 43.62050163576881 },
    @{ Row = 233; A = 44307; B = 0; C = 2; D = 43.62050163576881 }
)

# Carry the formatting of the last existing row (A229) down onto the new A cells
# before writing values, mirroring how the source data was produced (copy-down).
$ws.Range("A229").Copy() | Out-Null
$ws.Range("A230:A233").PasteSpecial(-4122) | Out-Null

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
}
